$d = $word.ActiveDocument

# The document currently has a "_GoBack" bookmark sitting in the middle of
# the answer paragraph for 2.10.3, wedged between two <w:tab/> runs. Word
# re-stamps this bookmark to mark the location of the most recent edit, and
# in the target revision that location is the first of the trailing empty
# paragraphs at the very end of the document (just after the last answer
# paragraph, before the final blank paragraphs / section break).

# Remove the existing _GoBack bookmark from its current location.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# The document ends with the last answer paragraph followed by four empty
# paragraphs and then the section break. Re-create the bookmark collapsed
# in the first of those four trailing empty paragraphs.
$n = $d.Paragraphs.Count
$target = $d.Paragraphs.Item($n - 3)
$d.Bookmarks.Add("_GoBack", $target.Range)
